$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.054.10"
$ws.Range("E2").Value = "  -15.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.219.42"
$ws.Range("E3").Value = "  -23.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.67"
$ws.Range("E5").Value = "  -20.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.53"
$ws.Range("E6").Value = "  -19.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.449"
$ws.Range("E8").Value = "  -18.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.162.68"
$ws.Range("E9").Value = "  -25.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "4.99"
$ws.Range("E10").Value = "  -17.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0851"
$ws.Range("E11").Value = "  -21.16%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.286"
$ws.Range("E12").Value = "  -20.56%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.118"
$ws.Range("E13").Value = "  -6.94%  "

$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "50.994.15"
$ws.Range("E14").Value = "  -15.92%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.08"
$ws.Range("E15").Value = "  -20.68%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000112"
$ws.Range("E16").Value = "  -20.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.208.50"
$ws.Range("E17").Value = "  -24.15%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "285.68"
$ws.Range("E18").Value = "  -20.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.74"
$ws.Range("E19").Value = "  -25.63%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("B21").Value = "LEO"
$ws.Range("C21").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.57"
$ws.Range("E21").Value = "  -1.68%  "

$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.41"
$ws.Range("E22").Value = "  -28.22%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.05"
$ws.Range("E23").Value = "  -23.90%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "51.93"
$ws.Range("E24").Value = "  -19.98%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.347"
$ws.Range("E25").Value = "  -23.72%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.134"
$ws.Range("E26").Value = "  -25.97%  "

$ws.Range("B27").Value = "USDe"
$ws.Range("C27").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.41%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.56"
$ws.Range("E28").Value = "  -16.76%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0643"
$ws.Range("E29").Value = "  -24.34%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.85"
$ws.Range("E30").Value = "  -6.68%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "16.31"
$ws.Range("E31").Value = "  -17.76%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("E32").Value = "  -26.44%  "

$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  -20.16%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.779"
$ws.Range("E34").Value = "  -22.49%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.30"
$ws.Range("E35").Value = "  -24.37%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.989"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.947"
$ws.Range("E37").Value = "  -21.37%  "

$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.18"
$ws.Range("E38").Value = "  -1.49%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.64"
$ws.Range("E39").Value = "  -19.23%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -19.98%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.16"
$ws.Range("E41").Value = "  -22.23%  "

$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0472"
$ws.Range("E42").Value = "  -19.54%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.815.85"
$ws.Range("E43").Value = "  -20.98%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.491"
$ws.Range("E44").Value = "  -24.42%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0197"
$ws.Range("E45").Value = "  -17.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0789"
$ws.Range("E46").Value = "  -14.45%  "

$ws.Range("B47").Value = "ZEEBU"
$ws.Range("C47").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.52"
$ws.Range("E47").Value = "  -6.65%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "14.71"
$ws.Range("E48").Value = "  -28.19%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.58"
$ws.Range("E49").Value = "  -28.37%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "14.31"
$ws.Range("E50").Value = "  -22.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.856"
$ws.Range("E51").Value = "  -15.69%  "
